# Update the two-digit-divided-by-one-digit practice problems.
# The worksheet has a single 20-row/5-column table; every 4th row
# (1, 5, 9, 13, 17) holds five division problems, the rows between
# them are blank spacer rows. Because a couple of the old problem
# strings repeat (e.g. "39÷3="), we address each problem cell
# positionally via Table.Cell(row, col) instead of a global
# Find/Replace, so each occurrence gets its own correct new value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)

$values = @(
    @("17÷6=", "31÷4=", "28÷5=", "55÷6=", "47÷6="),
    @("72÷5=", "60÷9=", "18÷9=", "87÷8=", "10÷3="),
    @("97÷5=", "41÷7=", "52÷9=", "95÷7=", "54÷9="),
    @("18÷3=", "95÷3=", "31÷5=", "21÷7=", "41÷2="),
    @("52÷2=", "61÷6=", "31÷5=", "37÷6=", "14÷3=")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    $rowValues = $values[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $t.Cell($row, $col).Range.Text = $rowValues[$col - 1]
    }
}
